$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Duplicated NIRs")
$lo = $ws.ListObjects.Item(1)

# Insert 3 new columns in the raw grid right before "protein_dry_basis" (col K)
$ws.Columns("K:M").Insert()
$ws.Range("K1").Value = "Pro13"
$ws.Range("L1").Value = "Oil13"
$ws.Range("M1").Value = "PO13"

# The inserted columns leave K2:M2 without cell records; copy an existing blank
# data cell into them so the row stays fully populated like the rest of row 2.
$ws.Range("A2").Copy($ws.Range("K2:M2"))

# Rebuild the table so its ListColumns metadata correctly reflects the new headers
# (Resize() alone leaves stale/duplicated column names in this runtime).
$lo.Unlist()
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:W2"), $true)
$lo2.Name = "Table3"
$lo2.TableStyle = "TableStyleLight9"
